# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" (fund holdings detail) right before
#    the "总计" (totals) sheet, formatted like the other quarterly sheets.
# 2) Insert a new first data row into "总计" summarizing the 2022-Q1 totals,
#    shifting the existing rows down and bumping their running index by 1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value into a cell as TEXT even when it looks numeric
# (keeps leading zeros / avoids numeric auto-conversion), while leaving the
# cell's style untouched (copies formatting back from a guaranteed-blank
# cell after the assignment so no stray NumberFormat-driven style lingers).
# ---------------------------------------------------------------------------
function Set-TextValue($cell, $text, $blankRef) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $blankRef.Copy()
    $cell.PasteSpecial(-4122)
}

# ===========================================================================
# Part 1: create the "2022-Q1" worksheet before "总计"
# ===========================================================================
$template = $wb.Worksheets.Item("2021-Q4")

# NOTE: `Worksheets.Add(beforeSheet)` inserts the new (blank) sheet at
# beforeSheet's position and the pre-fetched handle ends up tracking the
# *position* rather than the original "总计" sheet, so it must be looked
# up again (by name) after the insertion below, never reused beforehand.
$totalSheetRefForPosition = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheetRefForPosition)
$newSheet.Name = "2022-Q1"

# Copy header-row formatting (bold/centered/bordered header cells).
$template.Range("A1:H1").Copy()
$newSheet.Range("A1:H1").PasteSpecial(-4122)

# Copy a data-row formatting pattern down across all 8 data rows.
$template.Range("A2:H2").Copy()
$newSheet.Range("A2:H9").PasteSpecial(-4122)

$blankRef = $newSheet.Range("Z100")

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$rows = @(
    @{ idx = 0; code = "501030"; name = "汇添富中证环境治理指数（LOF）A"; scale = "6.61"; pos = "93.20"; pct = "2.12"; mv = "0.1401"; rank = 5 },
    @{ idx = 1; code = "501031"; name = "汇添富中证环境治理指数（LOF）C"; scale = "2.74"; pos = "93.20"; pct = "2.12"; mv = "0.0581"; rank = 5 },
    @{ idx = 2; code = "164908"; name = "交银施罗德中证环境治理指数（LOF）"; scale = "2.12"; pos = "93.72"; pct = "2.17"; mv = "0.0460"; rank = 5 },
    @{ idx = 3; code = "005632"; name = "鹏华量化先锋混合"; scale = "3.10"; pos = "92.91"; pct = "1.39"; mv = "0.0431"; rank = 9 },
    @{ idx = 4; code = "501219"; name = "华夏智胜先锋股票（LOF）A"; scale = "3.61"; pos = "94.50"; pct = "0.95"; mv = "0.0343"; rank = 10 },
    @{ idx = 5; code = "014198"; name = "华夏智胜先锋股票（LOF）C"; scale = "1.30"; pos = "94.50"; pct = "0.95"; mv = "0.0124"; rank = 10 },
    @{ idx = 6; code = "009263"; name = "华宝红利精选混合A"; scale = "0.46"; pos = "83.67"; pct = "0.97"; mv = "0.0045"; rank = 9 },
    @{ idx = 7; code = "010841"; name = "华宝红利精选混合C"; scale = "0.16"; pos = "83.67"; pct = "0.97"; mv = "0.0016"; rank = 9 }
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $row.idx
    Set-TextValue $newSheet.Cells.Item($r, 2) $row.code $blankRef
    $newSheet.Cells.Item($r, 3).Value = $row.name
    Set-TextValue $newSheet.Cells.Item($r, 4) $row.scale $blankRef
    Set-TextValue $newSheet.Cells.Item($r, 5) $row.pos $blankRef
    Set-TextValue $newSheet.Cells.Item($r, 6) $row.pct $blankRef
    Set-TextValue $newSheet.Cells.Item($r, 7) $row.mv $blankRef
    $newSheet.Cells.Item($r, 8).Value = $row.rank
    $r = $r + 1
}

# ===========================================================================
# Part 2: prepend a "2022-Q1" totals row into the "总计" sheet
# ===========================================================================
# Re-fetch "总计" by name now that the sheet list has been rearranged.
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# Bump the running index (column A) of all the pre-existing rows by 1,
# since they have all shifted down one row to make room for the new entry.
for ($rr = 7; $rr -ge 3; $rr--) {
    $cell = $totalSheet.Cells.Item($rr, 1)
    $cell.Value = $cell.Value2 + 1
}

# The freshly inserted row picked up some interpolated formatting; clear it
# on the non-index columns (they should carry no explicit style), then set
# the new values.
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 8
$totalSheet.Cells.Item(2, 4).Value = 0.34

# Restore column-A's bold/centered/bordered style on the new row (copy it
# from the row right below, which still carries the original styling).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

Write-Host "2022-Q1 sheet added and 总计 sheet updated"
